$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1436130007558579
$ws.Range("C2").Value = 0.1332829428067523
$ws.Range("D2").Value = 0.1421012849584278
$ws.Range("E2").Value = 0.1428571428571428
$ws.Range("F2").Value = 0.1337868480725624
$ws.Range("G2").Value = 0.1330309901738473
$ws.Range("I2").Value = 0.2151675485008818
$ws.Range("J2").Value = 0.2204585537918871
$ws.Range("K2").Value = 0.2078609221466364
$ws.Range("L2").Value = 0.199798437893676
$ws.Range("M2").Value = 0.1987906273620559
$ws.Range("B3").Value = 0.1295036533131771
$ws.Range("C3").Value = 0.1352985638699924
$ws.Range("D3").Value = 0.1405895691609977
$ws.Range("E3").Value = 0.1443688586545729
$ws.Range("F3").Value = 0.1468883849836231
$ws.Range("G3").Value = 0.1398337112622827
$ws.Range("I3").Value = 0.2093726379440665
$ws.Range("J3").Value = 0.2244897959183673
$ws.Range("K3").Value = 0.2106324011085916
$ws.Range("L3").Value = 0.2136558327034518
$ws.Range("M3").Value = 0.1922398589065256
$ws.Range("B4").Value = 0.1201814058956916
$ws.Range("C4").Value = 0.1209372637944066
$ws.Range("D4").Value = 0.1186696900982615
$ws.Range("E4").Value = 0.1211892164273117
$ws.Range("F4").Value = 0.1234567901234568
$ws.Range("G4").Value = 0.1211892164273117
$ws.Range("H4").Value = 0.1680524061476442
$ws.Range("I4").Value = 0.1710758377425044
$ws.Range("K4").Value = 0.1657848324514991
$ws.Range("L4").Value = 0.162761400856639
$ws.Range("M4").Value = 0.1569664902998237
$ws.Range("B5").Value = 0.1070798689846309
$ws.Range("C5").Value = 0.09397833207357018
$ws.Range("D5").Value = 0.109599395313681
$ws.Range("E5").Value = 0.1055681531872008
$ws.Range("F5").Value = 0.1126228269085412
$ws.Range("G5").Value = 0.1070798689846309
$ws.Range("H5").Value = 0.1574703955656337
$ws.Range("I5").Value = 0.1592340639959688
$ws.Range("J5").Value = 0.163013353489544
$ws.Range("K5").Value = 0.1544469639707735
$ws.Range("M5").Value = 0.145628621819098
$ws.Range("B6").Value = 0.1194255479969766
$ws.Range("D6").Value = 0.1121189216427312
$ws.Range("E6").Value = 0.1138825900730663
$ws.Range("F6").Value = 0.1156462585034014
$ws.Range("G6").Value = 0.1012849584278156
$ws.Range("H6").Value = 0.1652809271856891
$ws.Range("I6").Value = 0.1667926429831192
$ws.Range("J6").Value = 0.1617535903250189
$ws.Range("K6").Value = 0.1599899218946838
$ws.Range("L6").Value = 0.1604938271604938
$ws.Range("M6").Value = 0.1494079113126732
$ws.Range("H7").Value = 0.1347946586041824
$ws.Range("I7").Value = 0.1350466112370874
$ws.Range("J7").Value = 0.145376669186193
$ws.Range("K7").Value = 0.1257243638196019
$ws.Range("B8").Value = 0.1340388007054674
$ws.Range("C8").Value = 0.1335348954396574
$ws.Range("D8").Value = 0.1335348954396574
$ws.Range("E8").Value = 0.1345427059712774
$ws.Range("F8").Value = 0.1365583270345175
$ws.Range("G8").Value = 0.1332829428067523
$ws.Range("H8").Value = 0.1899722852103805
$ws.Range("I8").Value = 0.1841773746535651
$ws.Range("J8").Value = 0.1723356009070295
$ws.Range("K8").Value = 0.163517258755354
$ws.Range("L8").Value = 0.1589821113630638
$ws.Range("M8").Value = 0.1594860166288738
$ws.Range("B9").Value = 0.1045603426555807
$ws.Range("C9").Value = 0.1292517006802721
$ws.Range("D9").Value = 0.1289997480473671
$ws.Range("E9").Value = 0.1307634164777022
$ws.Range("F9").Value = 0.1312673217435122
$ws.Range("G9").Value = 0.1257243638196019
$ws.Range("I9").Value = 0.1849332325522802
$ws.Range("J9").Value = 0.1791383219954649
$ws.Range("K9").Value = 0.1690602166792643
$ws.Range("L9").Value = 0.1672965482489292
$ws.Range("M9").Value = 0.164273116654069
$ws.Range("C10").Value = 0.1111111111111111
$ws.Range("D10").Value = 0.1063240110859159
$ws.Range("E10").Value = 0.1075837742504409
$ws.Range("F10").Value = 0.1128747795414462
$ws.Range("H10").Value = 0.1589821113630638
$ws.Range("I10").Value = 0.1516754850088183
$ws.Range("J10").Value = 0.1441169060216679
$ws.Range("K10").Value = 0.145124716553288
$ws.Range("L10").Value = 0.1365583270345175
$ws.Range("M10").Value = 0.1340388007054674
$ws.Range("C11").Value = 0.08339632149155958
$ws.Range("E11").Value = 0.09599395313681028
$ws.Range("F11").Value = 0.09221466364323506
$ws.Range("G11").Value = 0.08667170571932477
$ws.Range("H11").Value = 0.1506676744771983
$ws.Range("I11").Value = 0.1489040060468632
$ws.Range("J11").Value = 0.145376669186193
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.1408415217939027
$ws.Range("M11").Value = 0.1388259007306626
$ws.Range("I12").Value = 0.1370622323003275
$ws.Range("J12").Value = 0.1257243638196019
$ws.Range("K12").Value = 0.1330309901738473
$ws.Range("L12").Value = 0.1365583270345175
$ws.Range("M12").Value = 0.1164021164021164
$ws.Range("H13").Value = 0.1224489795918367
$ws.Range("I13").Value = 0.1209372637944066
$ws.Range("J13").Value = 0.1101033005794911
$ws.Range("K13").Value = 0.1131267321743512
$ws.Range("L13").Value = 0.09549004787100025
$ws.Range("C14").Value = 0.1128747795414462
$ws.Range("D14").Value = 0.1073318216175359
$ws.Range("E14").Value = 0.1106072058453011
$ws.Range("F14").Value = 0.1108591584782061
$ws.Range("G14").Value = 0.1035525321239607
$ws.Range("H14").Value = 0.1423532375913328
$ws.Range("I14").Value = 0.1385739480977576
$ws.Range("J14").Value = 0.1403376165280927
$ws.Range("K14").Value = 0.1375661375661376
$ws.Range("L14").Value = 0.145124716553288
$ws.Range("M14").Value = 0.1365583270345175
$ws.Range("C15").Value = 0.08390022675736961
$ws.Range("D15").Value = 0.09045099521289998
$ws.Range("E15").Value = 0.09120685311161501
$ws.Range("F15").Value = 0.09221466364323506
$ws.Range("G15").Value = 0.09322247417485512
$ws.Range("H15").Value = 0.1405895691609977
$ws.Range("I15").Value = 0.1400856638951877
$ws.Range("J15").Value = 0.1297556059460821
$ws.Range("K15").Value = 0.1360544217687075
$ws.Range("M15").Value = 0.1385739480977576
$ws.Range("D16").Value = 0.07180650037792895
$ws.Range("E16").Value = 0.07432602670697909
$ws.Range("F16").Value = 0.08616780045351474
$ws.Range("G16").Value = 0.06223230032753842
$ws.Range("H16").Value = 0.1161501637692114
$ws.Range("I16").Value = 0.1126228269085412
$ws.Range("J16").Value = 0.1146384479717813
$ws.Range("K16").Value = 0.109599395313681
$ws.Range("L16").Value = 0.1103552532123961
$ws.Range("M16").Value = 0.1121189216427312
$ws.Range("D17").Value = 0
$ws.Range("H17").Value = 0.1123708742756362
$ws.Range("I17").Value = 0.109599395313681
$ws.Range("J17").Value = 0.1058201058201058
$ws.Range("K17").Value = 0.1025447215923406
$ws.Range("L17").Value = 0.1038044847568657
$ws.Range("M17").Value = 0.1068279163517259
$ws.Range("H18").Value = 0.09322247417485512
$ws.Range("I18").Value = 0.09523809523809523
$ws.Range("J18").Value = 0.09926933736457547
$ws.Range("K18").Value = 0.09322247417485512
$ws.Range("L18").Value = 0.09624590576971528
$ws.Range("M18").Value = 0.09624590576971528
$ws.Range("H19").Value = 0.06475182665658856
$ws.Range("I19").Value = 0.06349206349206349
$ws.Range("J19").Value = 0.06500377928949357
$ws.Range("K19").Value = 0.06752330561854371
$ws.Range("L19").Value = 0.07583774250440917
$ws.Range("M19").Value = 0.0818846056941295
$ws.Range("C20").Value = 0.07936507936507936
$ws.Range("D20").Value = 0.1131267321743512
$ws.Range("E20").Value = 0.1063240110859159
$ws.Range("F20").Value = 0.1201814058956916
$ws.Range("G20").Value = 0.1164021164021164
$ws.Range("H20").Value = 0.1884605694129504
$ws.Range("I20").Value = 0.1909800957420005
$ws.Range("J20").Value = 0.1950113378684807
$ws.Range("K20").Value = 0.201058201058201
$ws.Range("L20").Value = 0.2078609221466364
$ws.Range("M20").Value = 0.1980347694633409
$ws.Range("C21").Value = 0.06147644242882338
$ws.Range("D21").Value = 0.1027966742252457
$ws.Range("E21").Value = 0.1133786848072562
$ws.Range("F21").Value = 0.1365583270345175
$ws.Range("G21").Value = 0.1534391534391534
$ws.Range("H21").Value = 0.2348198538674729
$ws.Range("I21").Value = 0.218442932728647
$ws.Range("J21").Value = 0.217687074829932
$ws.Range("K21").Value = 0.2101284958427815
$ws.Range("L21").Value = 0.2338120433358528
$ws.Range("M21").Value = 0.2413706223230033
$ws.Range("B22").Value = 0.08112874779541446
$ws.Range("C22").Value = 0.05467372134038801
$ws.Range("D22").Value = 0.08868732678256487
$ws.Range("E22").Value = 0.1048122952884858
$ws.Range("F22").Value = 0.1161501637692114
$ws.Range("G22").Value = 0.1169060216679264
$ws.Range("H22").Value = 0.1897203325774754
$ws.Range("I22").Value = 0.181657848324515
$ws.Range("J22").Value = 0.1516754850088183
$ws.Range("K22").Value = 0.127991937515747
$ws.Range("L22").Value = 0.1249685059208869
$ws.Range("M22").Value = 0.2030738221214412
$ws.Range("C23").Value = 0.0564373897707231
$ws.Range("D23").Value = 0.05668934240362812
$ws.Range("E23").Value = 0.05769715293524817
$ws.Range("F23").Value = 0.06046863189720332
$ws.Range("H23").Value = 0.1698160745779793
$ws.Range("I23").Value = 0.1519274376417233
$ws.Range("J23").Value = 0.1748551272360796
$ws.Range("K23").Value = 0.1242126480221718
$ws.Range("L23").Value = 0.1418493323255228
$ws.Range("M23").Value = 0.2474174855127236
